# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" before the "总计" sheet, populated
#    with the fund-holding detail rows for the new quarter.
# 2. Update the "总计" (totals) summary sheet with a new leading row for
#    2022-Q1, shifting the previously existing rows down by one.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")

# --- Create the new "2022-Q1" sheet just before "总计" ---------------------
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# NOTE: after Add() runs, the $total variable ends up referring to the same
# sheet object as $newSheet, so we must re-acquire a fresh reference to the
# "总计" sheet by name before touching it again below.
$total = $wb.Worksheets.Item("总计")

# --- Header row --------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data rows (columns B-G are stored as text, matching the source data) --
$textRange = $newSheet.Range("B2:G3")
$textRange.NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "006601"
$newSheet.Range("C2").Value = "国融融泰灵活配置混合A"
$newSheet.Range("D2").Value = "0.04"
$newSheet.Range("E2").Value = "47.44"
$newSheet.Range("F2").Value = "2.77"
$newSheet.Range("G2").Value = "0.0011"
$newSheet.Range("H2").Value = 7

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "006602"
$newSheet.Range("C3").Value = "国融融泰灵活配置混合C"
$newSheet.Range("D3").Value = "0.01"
$newSheet.Range("E3").Value = "47.44"
$newSheet.Range("F3").Value = "2.77"
$newSheet.Range("G3").Value = "0.0003"
$newSheet.Range("H3").Value = 7

$indexRange = $newSheet.Range("A2:A3")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

# --- Update the "总计" sheet: push existing rows down and insert the new
#     2022-Q1 summary row at the top of the data (row 2). ------------------
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.06

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0

$total.Range("A4").Font.Bold = $true
$total.Range("A4").HorizontalAlignment = -4108
$total.Range("A4").VerticalAlignment = -4160
$total.Range("A4").Borders.LineStyle = 1

Write-Output "applied 2022-Q1 edit"
